# "updated parts list, including pricing"
#
# case_1 (sheet1): the KY-016 indicator LED row (row 17) switches supplier
# from Amazon to Banggood, with a new (lower) unit price and an updated
# notes/alternatives string warning about COVID shipping delays. The old
# "large kit of sensors" note text is no longer referenced by any cell, so
# the shared-string table drops it and appends the two brand-new strings
# ("Banggood" + the COVID note) - the engine re-indexes every other shared
# string reference automatically. The SUM total in C18 recalculates on its
# own once C17 changes.
#
# The hyperlink that used to sit on the KY-016 row (an Amazon listing for a
# sensor kit) moves to the jumper-wire row (A5) instead, while A5's old
# AdaFruit link moves up to the KY-016 row (A17) - i.e. the two hyperlink
# targets trade places.
#
# Both sheets also record a different active cell/tab: case_1 becomes the
# selected tab (active cell D5), and case_10 stops being the selected tab
# (active cell A13).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("case_1")
$ws2 = $wb.Worksheets.Item("case_10")

# --- case_1 row 17: KY-016 indicator LED now sourced from Banggood ---
$ws1.Range("B17").Value2 = "Banggood"
$ws1.Range("C17").Value2 = 4.27
$ws1.Range("D17").Value2 = "ALLOW 3+ WEEKS TO ARRIVE FROM CHINA. Due to COVID, it's hard to find these in the US. You can get it quickly from Amazon B07KJYR8K1, but costs `$18."

# --- swap the A5 / A17 hyperlink targets (keep each cell's own formatting) ---
$hlA17 = $null
$hlA5  = $null
foreach ($hl in $ws1.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$17") { $hlA17 = $hl }
    if ($addr -eq "`$A`$5")  { $hlA5  = $hl }
}
$addrForA17 = $hlA17.Address
$addrForA5  = $hlA5.Address

$fontSizeA5  = $ws1.Range("A5").Font.Size
$fontSizeA17 = $ws1.Range("A17").Font.Size

$hlA5.Delete()
$hlA17.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A5"), $addrForA17) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A17"), $addrForA5) | Out-Null

# re-adding a hyperlink resets the cell to Excel's default Hyperlink style;
# put the original (smaller) font size back so the row keeps its look.
$ws1.Range("A5").Font.Size = $fontSizeA5
$ws1.Range("A17").Font.Size = $fontSizeA17

# --- active tab / selection: case_1 becomes the shown sheet ---
$ws1.Activate()
$ws1.Range("D5").Select()
$ws2.Activate()
$ws2.Range("A13").Select()
$ws1.Activate()
